# AIP-1584 AIP-1586 - Updated test data for some of the cablings
#
# 16-channel cabling dataset: trims the DeviceInfo calibration/IP pair,
# drops the two extra (17th/18th) channel rows from the Cabling table,
# and zeroes out the now-unused 5th-feeder / DSP-channel slots in the
# BusbarFeederMap and DSPChannelMap sheets.

$wb = $excel.ActiveWorkbook

# --- DeviceInfo: point at the new device / calibration file ---------------
$wsDevice = $wb.Worksheets.Item("DeviceInfo")
$wsDevice.Range("A2").Value = "10.75.58.66"
$wsDevice.Range("B2").Value = "6U_10I.cal"

# --- Cabling: this is a 16-channel config, so Channel[16]/Channel[17] ------
# (rows 18-19, the 17th and 18th table rows) no longer apply - clear them.
$wsCabling = $wb.Worksheets.Item("Cabling")
$wsCabling.Range("A18:J19").ClearContents() | Out-Null

# --- BusbarFeederMap: Feeder 4 only used 2 channels before; now unused -----
$wsBusbar = $wb.Worksheets.Item("BusbarFeederMap")
$wsBusbar.Range("B36").Value = 0   # Feeder_4_Channel_0
$wsBusbar.Range("B37").Value = 0   # Feeder_4_Channel_1
$wsBusbar.Range("B40").Value = 0   # Feeder_4_Channel_Count

# --- DSPChannelMap: DSP2 channels 6/7 no longer mapped ---------------------
$wsDsp = $wb.Worksheets.Item("DSPChannelMap")
$wsDsp.Range("D8").Value = 0
$wsDsp.Range("D9").Value = 0

# --- View/selection housekeeping: select-all on every sheet, leaving ------
# DeviceInfo as the active tab (matches the saved workbook's view state).
foreach ($name in @("DeviceInfo", "Cabling", "DSPFeederMap", "BusbarFeederMap", "DSPChannelMap")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate() | Out-Null
    $ws.Cells.Select() | Out-Null
}
$wsDevice.Activate() | Out-Null
